# Update gh-pages to output generated at 456a3b4
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 352
$ws1.Range("F4").Value = 10634
$ws1.Range("F6").Value = 967
$ws1.Range("F7").Value = 121
$ws1.Range("F8").Value = 1312
$ws1.Range("F9").Value = 8233
$ws1.Range("G10").Value = 89
$ws1.Range("F13").Value = 214
$ws1.Range("F14").Value = 135
$ws1.Range("F15").Value = 3269
$ws1.Range("F18").Value = 753
$ws1.Range("F23").Value = 1718

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 352
$ws4.Range("F4").Value = 10634
$ws4.Range("F6").Value = 967
$ws4.Range("F7").Value = 121
$ws4.Range("F8").Value = 1312
$ws4.Range("F9").Value = 8233
$ws4.Range("G10").Value = 89
$ws4.Range("F13").Value = 214
$ws4.Range("F14").Value = 135
$ws4.Range("F15").Value = 3269
$ws4.Range("F18").Value = 753
$ws4.Range("F23").Value = 1718
